$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add headers I1 "I0" and J1 "IF", copying the header style (s="1") from H1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the I and J column data for rows 2-61
$iVals = @(6,8,8,7,7,8,7,6,7,9,6,7,7,7,8,5,7,7,8,9,6,8,8,8,6,6,8,4,5,6,7,8,9,8,10,6,9,7,8,8,7,9,8,9,5,8,9,9,4,6,6,10,6,6,6,6,5,5,5,5)
$jVals = @(6,8,8,7,7,8,7,7,7,9,7,7,7,7,8,5,7,7,8,9,6,8,8,8,6,7,8,4,5,6,7,8,9,8,10,6,9,7,8,8,7,9,8,9,6,8,9,9,6,7,6,10,6,6,6,6,5,5,5,5)
for ($r = 0; $r -lt $iVals.Length; $r++) {
    $ws.Cells.Item($r + 2, 9).Value = $iVals[$r]
    $ws.Cells.Item($r + 2, 10).Value = $jVals[$r]
}
